$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Test Cases")

# "Running suites a and b" - flip Runmode from "N" to "Y" for every test
# case (rows 2-16) so the whole suite executes, matching row 17 which was
# already set to run.
$ws.Range("C2:C16").Value = "Y"

# Match the formatting of the already-"Y" row (C17 picks up the same
# fill-less cell style as the rest of the column).
$ws.Range("C16").Copy()
[void]$ws.Range("C17").PasteSpecial(-4122)

# Leave the selection where the edit was made.
[void]$ws.Range("C3:C17").Select()
